# COMP6015 Report Plan - "Minor Report Plan Update"
#
# 1) "Section 6 - References" (plain hyphen, single run) becomes three
#    runs: "Section 6 " + an en dash "\u2013" + " References".
# 2) A new "Next Steps" Heading 1 section is appended after it, followed
#    by a blank paragraph and two body paragraphs ("Research main
#    security features" / "Make Notes"), ahead of the document's
#    existing trailing blank paragraph.

$d = $word.ActiveDocument

# Locate the "Section 6 - References" paragraph without assuming a fixed
# paragraph index.
$target = $d.Content
$found = $target.Find.Execute("Section 6 - References", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find 'Section 6 - References' paragraph"
}

$enDash = [string][char]0x2013

# Flat-OPC WordprocessingML fragment: the first <w:p> reproduces the
# "Section 6 ... References" line split across three runs (so the dash
# is its own run), then four brand-new paragraphs follow it.
$xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Section 6 </w:t></w:r><w:r><w:t>$enDash</w:t></w:r><w:r><w:t xml:space="preserve"> References</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t>Next Steps</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>Research main security features</w:t></w:r></w:p><w:p><w:r><w:t>Make Notes</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$target.InsertXML($xml, "Replace") | Out-Null
